$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values in row 6
$ws.Range("E6").Value = 6
$ws.Range("G6").Value = -3
$ws.Range("H6").Value = 13

# Update the selected cell to D6
$ws.Range("D6").Select()
